$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("A29").Value = "boolean"
$ws.Range("B29").Value = $true
$ws.Range("C29").Value = $true
$ws.Range("D29").Value = $false

$ws.Range("C39").Select()
